# "aggiornamento fino a 28 luglio" - append new daily rows (302-328) to the
# data table, continuing the date series (Excel serial dates) and filling
# the "nuovi pos." / "somma mobile 7gg." / "somma mobile 7gg. per 100mila
# abitanti" columns with 0, matching the existing rows' layout/format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 302
$endRow = 328
$startSerial = 44376

# Reference cell whose style (date number format, border, alignment) the new
# date cells in column A should inherit - same style used by every existing
# row in column A.
$styleSource = $ws.Range("A301")

for ($i = 0; $i -le ($endRow - $startRow); $i++) {
    $row = $startRow + $i
    $serial = $startSerial + $i

    # Copy formatting (incl. number format/border/alignment) onto the new
    # date cell before writing its value.
    $styleSource.Copy() | Out-Null
    $ws.Cells.Item($row, 1).PasteSpecial(-4122) | Out-Null

    $ws.Cells.Item($row, 1).Value = $serial
    $ws.Cells.Item($row, 2).Value = 0
    $ws.Cells.Item($row, 3).Value = 0
    $ws.Cells.Item($row, 4).Value = 0
}

$excel.CutCopyMode = $false
